$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1, 1).End(4).Row

$colC = $ws.Range($ws.Cells.Item(1, 3), $ws.Cells.Item($lastRow, 3))
$colD = $ws.Range($ws.Cells.Item(1, 4), $ws.Cells.Item($lastRow, 4))

$valsC = $colC.Value2
$valsD = $colD.Value2

$colC.Value = $valsD
$colD.Value = $valsC
